$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: newly scraped job record appended below the existing header row.
# (Summary/D2 was scraped as an empty string, so it is intentionally left blank.)
$ws.Range("A2").Value = "Mobile Application Developer - Flutter"
$ws.Range("B2").Value = "Future Focus Infotech"
$ws.Range("C2").Value = "Ajman, UAE"
$ws.Range("E2").Value = "13 days ago"
$ws.Range("F2").Value = "https://www.gulftalent.com/uae/jobs/mobile-application-developer-flutter-375078"
